$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# Query text blocks (order chosen so the workbook shared-string table
# fills in the same sequence Excel produced them in)
$filesQuery = 'MATCH (f:file)-->(s:study)
OPTIONAL MATCH (samp:sample)<--(f)
OPTIONAL MATCH (samp)-->(p:participant)
OPTIONAL MATCH (f)<--(g:genomic_info)
OPTIONAL MATCH (p)<--(diag:diagnosis)
WITH s, p, samp, f, g, diag
WHERE g.library_strategy in  [''Targeted-Capture'']
WITH DISTINCT f, s, p, samp
RETURN
    coalesce(f.file_name, '''') as `File Name`,
    coalesce(s.study_name,'''') as `Study Name`,
    coalesce(s.phs_accession,'''') as `Accession`,
    coalesce(p.participant_id, '''') as `Participant ID`,
    coalesce(samp.sample_id, '''') as `Sample ID`,
    coalesce(f.file_type, '''') as `File Type`
ORDER BY f.file_name limit 100'
$samplesQuery = 'MATCH (samp:sample)-->(p:participant)-->(s:study)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
OPTIONAL MATCH (p)<--(diag:diagnosis)
WITH s, p, samp, f, g, diag
WHERE g.library_strategy in  [''Targeted-Capture'']
WITH DISTINCT s, p, samp
RETURN
    coalesce(samp.sample_id, '''') as `Sample ID`,
    coalesce(p.participant_id,'''') as `Participant ID`,
    coalesce(s.study_name, '''') as `Study Name`,
    coalesce(s.phs_accession,'''') as `Accession`,
    coalesce(samp.sample_tumor_status,'''') as `Tumor`,
    coalesce(samp.sample_type,'''') as `Analyte Type`
ORDER BY samp.sample_id limit 100'
$participantsQuery = 'MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH s, p, samp, f, g, diag
WHERE g.library_strategy in [''Targeted-Capture'']
WITH p
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp
RETURN 
coalesce(p.participant_id,'''') as `Participant ID`,
coalesce(s.study_name, '''') as `Study Name`,
coalesce(s.phs_accession,'''') as `Accession`,
coalesce(p.gender,'''') as `Gender`,
coalesce(apoc.text.join(samp, '',''), '''') as `Samples`
ORDER BY p.participant_id limit 100'
$statQuery = 'CALL{
    MATCH (p:participant)-->(s:study)
    OPTIONAL MATCH (samp:sample)-->(p)
    OPTIONAL MATCH (samp)<--(f:file)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag
    WHERE g.library_strategy in  [''Targeted-Capture'']
    RETURN 
        count(distinct p) AS num_participants
}
WITH num_participants
CALL {
    MATCH (samp:sample)-->(p:participant)-->(s)
    OPTIONAL MATCH (samp)<--(f:file)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag
    WHERE g.library_strategy in  [''Targeted-Capture'']
    RETURN 
        count(distinct samp) AS num_samples
}
WITH num_participants, num_samples
CALL {
    MATCH (f:file)-->(s:study)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (samp:sample)<--(f)
    OPTIONAL MATCH (p:participant)<--(samp)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag
    WHERE g.library_strategy in  [''Targeted-Capture'']
    RETURN 
        count(distinct s) AS num_studies,
        count(distinct f) AS num_files
}
RETURN 
    num_studies AS Studies,
    num_participants AS Participants,
    num_samples AS Samples,
    num_files AS `Files`'

# --- Row 2: CasesTab -> ParticipantsTab ---
$ws.Range("A2").Value = "ParticipantsTab"
$ws.Range("C2").Value = $statQuery
$ws.Range("B4").Value = $filesQuery
$ws.Range("B3").Value = $samplesQuery
$ws.Range("B2").Value = $participantsQuery

# --- Row 3: SamplesTab (query rewritten) ---
$ws.Range("A3").Value = "SamplesTab"
$ws.Range("C3").Value = $statQuery

# --- Row 4: FilesTab (query rewritten) ---
$ws.Range("A4").Value = "FilesTab"
$ws.Range("C4").Value = $statQuery

# --- Row heights (grew to fit the longer rewritten Cypher queries) ---
$ws.Rows.Item(2).RowHeight = 386.25
$ws.Rows.Item(3).RowHeight = 330
$ws.Rows.Item(4).RowHeight = 369

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 22.666666666666668
$ws.Columns.Item(2).ColumnWidth = 95.83333333333333
$ws.Columns.Item(3).ColumnWidth = 74.0
$ws.Columns.Item(4).ColumnWidth = 78.0
$ws.Columns.Item(5).ColumnWidth = 76.33333333333333

# --- Selection / view ---
$ws.Range("D2").Select()
